# Add text content stripping.
# Update the SnippetID (column H) values in the voice-lines table with new
# randomly-regenerated 4-character snippet identifiers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "w7zp"
$ws.Range("H3").Value = "w7zp"
$ws.Range("H4").Value = "D5hS"
$ws.Range("H5").Value = "WSyU"
$ws.Range("H6").Value = "aaR0"
$ws.Range("H7").Value = "YDaQ"
$ws.Range("H8").Value = "blNR"
$ws.Range("H9").Value = "GX53"
$ws.Range("H10").Value = "nvdw"
$ws.Range("H11").Value = "SlpG"
$ws.Range("H12").Value = "sd2u"
$ws.Range("H13").Value = "LTD9"
$ws.Range("H14").Value = "wEui"
$ws.Range("H15").Value = "GVje"
$ws.Range("H16").Value = "7inI"
$ws.Range("H17").Value = "y2Bv"
$ws.Range("H18").Value = "y2Bv"
$ws.Range("H19").Value = "y2Bv"
$ws.Range("H20").Value = "y2Bv"
$ws.Range("H21").Value = "y2Bv"
$ws.Range("H22").Value = "hUaM"
$ws.Range("H23").Value = "lAkH"
$ws.Range("H24").Value = "PSMq"
$ws.Range("H25").Value = "xSfR"
$ws.Range("H26").Value = "U06q"
$ws.Range("H27").Value = "U06q"
$ws.Range("H28").Value = "EhGW"
$ws.Range("H29").Value = "ylr1"
